$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("catalogo")
$ws.Rows.Item(34).Delete()

$datos = $wb.Worksheets.Item("datos")
$datos.Range("C9").ClearContents()
